$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to the used range (A1:B5) - matches the new
# cellXfs entry (numFmtId 49 = "@") referenced by every cell's s="1".
$ws.Range("A1:B5").NumberFormat = "@"

# Update the customer data: B3/B4/B5 get new values. Because the Text
# format was applied first, these are stored as text (shared strings),
# matching kkjk585/hhjk44/55ggg -> 1345/kjhgff/78554 in the diff.
$ws.Range("B3").Value = "1345"
$ws.Range("B4").Value = "kjhgff"
$ws.Range("B5").Value = "78554"

# Switch the sheet's print orientation to portrait, which adds the
# <pageSetup orientation="portrait"/> element seen in the diff.
$ws.PageSetup.Orientation = 1
